# "the start of applications" - add two new dataset rows to the "data" sheet
# and update the existing transactions row's source file name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 3: the transactions export file name was refreshed.
$ws.Cells.Item(3, 6).Value = "transactions_11122022.xlsx"

# Row 4: new "calories" dataset.
$ws.Cells.Item(4, 1).Value = "d67af7eb-76b5-4943-9d4a-00b152c16171"
# Force the date-like string to be stored as literal text (not an actual
# date) the same way the existing rows store their "created" column.
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "11/04/2022"
$ws.Cells.Item(4, 4).Value = "calories"
$ws.Cells.Item(4, 6).Value = "calories_11042022.xlsx"
$ws.Cells.Item(4, 7).Value = "active"

# Row 5: new "calendar" dataset.
$ws.Cells.Item(5, 1).Value = "b27aaea9-db0f-4987-9771-700e2a6a43f5"
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "11/06/2022"
$ws.Cells.Item(5, 4).Value = "calendar"
$ws.Cells.Item(5, 6).Value = "calendar.xlsx"
$ws.Cells.Item(5, 7).Value = "active"

# Clone the formatting of row 2 (A is bold/bordered, B-H are plain) onto the
# two new rows so every new cell ends up with exactly the same cell style as
# the matching column in the existing data rows (and so that the otherwise
# unused C/E/H columns still get an (empty) cell written out).
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(2, $col).Copy()
    $ws.Cells.Item(4, $col).PasteSpecial(-4122)
    $ws.Cells.Item(2, $col).Copy()
    $ws.Cells.Item(5, $col).PasteSpecial(-4122)
}
